$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q3" right before the existing
#    "2022-Q2" sheet (currently the 2nd sheet in the workbook), and fill it
#    with the new quarter's fund-holding detail table.
# ---------------------------------------------------------------------------
$placeholder = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($placeholder)
$newSheet.Name = "2022-Q3"

# After insertion the worksheet that used to be at position 2 ("2022-Q2")
# has been pushed to position 3 - fetch it fresh so we can copy its header
# formatting onto the new sheet.
$oldQ2 = $wb.Worksheets.Item(3)

$oldQ2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$oldQ2.Range("A2:A3").Copy($newSheet.Range("A2:A3"))

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1

# These columns hold numeric-looking text (fund codes with leading zeros,
# decimal figures stored as text) - force text formatting so values are not
# silently coerced into numbers.
$newSheet.Range("B2:C3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "004405"
$newSheet.Range("C2").Value = "国寿安保稳寿混合A"
$newSheet.Range("D2").Value = "4.69"
$newSheet.Range("E2").Value = "24.56"
$newSheet.Range("F2").Value = "0.83"
$newSheet.Range("G2").Value = "0.0389"
$newSheet.Range("H2").Value = 10

$newSheet.Range("B3").Value = "004406"
$newSheet.Range("C3").Value = "国寿安保稳寿混合C"
$newSheet.Range("D3").Value = "0.62"
$newSheet.Range("E3").Value = "24.56"
$newSheet.Range("F3").Value = "0.83"
$newSheet.Range("G3").Value = "0.0051"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add a new row for 2022-Q3 holdings
#    and push the older quarters down, appending 2020-Q4 as the new last row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Copy formatting from the last existing data row down onto the newly
# appended row so it keeps the same style as the others (e.g. column A's
# index styling).
$summary.Range("A8").Copy($summary.Range("A9"))

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.04

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 15
$summary.Range("D3").Value = 1.52

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 41
$summary.Range("D4").Value = 10.31

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 47
$summary.Range("D5").Value = 11.15

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 24
$summary.Range("D6").Value = 9.03

$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 87
$summary.Range("D7").Value = 21.48

$summary.Range("B8").Value = "2021-Q1"
$summary.Range("C8").Value = 18
$summary.Range("D8").Value = 2.98

$summary.Range("A9").Value = 7
$summary.Range("B9").Value = "2020-Q4"
$summary.Range("C9").Value = 12
$summary.Range("D9").Value = 1.46
